$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "T.6"

$newSheet.Range("D5").Value = "Rank"
$newSheet.Range("E5").Value = "Framework"
$newSheet.Range("F5").Value = "Setup and Configuration"
$newSheet.Range("G5").Value = "Examples and Tutorials"
$newSheet.Range("H5").Value = "Custom Algorithm Implementation"
$newSheet.Range("I5").Value = "Adaptability to Various Use Cases"
$newSheet.Range("J5").Value = "Industry Adoption"
$newSheet.Range("K5").Value = "Average Score"

$newSheet.Range("D6").Value = "1st"
$newSheet.Range("D7").Value = "2nd"
$newSheet.Range("D8").Value = "3rd"
$newSheet.Range("D9").Value = "4th"
$newSheet.Range("D10").Value = "5th"

$newSheet.Range("F4").Value = "Ease of use"
$newSheet.Range("H4").Value = "Flexibility and customisability"
$newSheet.Range("J4").Value = "Real-world applicability"

$newSheet.Range("E6").Value = "FedML"
$newSheet.Range("F6").Value = 8
$newSheet.Range("G6").Value = 8
$newSheet.Range("H6").Value = 7
$newSheet.Range("I6").Value = 8
$newSheet.Range("J6").Value = 7
$newSheet.Range("K6").Value = 7.6

$newSheet.Range("E7").Value = "Flower"
$newSheet.Range("F7").Value = 7
$newSheet.Range("G7").Value = 8
$newSheet.Range("H7").Value = 8
$newSheet.Range("I7").Value = 8
$newSheet.Range("J7").Value = 6
$newSheet.Range("K7").Value = 7.4

$newSheet.Range("E8").Value = "TFF"
$newSheet.Range("F8").Value = 6
$newSheet.Range("G8").Value = 7
$newSheet.Range("H8").Value = 9
$newSheet.Range("I8").Value = 7
$newSheet.Range("J8").Value = 7
$newSheet.Range("K8").Value = 7.2

$newSheet.Range("E9").Value = "FATE"
$newSheet.Range("F9").Value = 6
$newSheet.Range("G9").Value = 6
$newSheet.Range("H9").Value = 7
$newSheet.Range("I9").Value = 6
$newSheet.Range("J9").Value = 6
$newSheet.Range("K9").Value = 6.2

$newSheet.Range("E10").Value = "PySyft"
$newSheet.Range("F10").Value = 5
$newSheet.Range("G10").Value = 5
$newSheet.Range("H10").Value = 6
$newSheet.Range("I10").Value = 5
$newSheet.Range("J10").Value = 4
$newSheet.Range("K10").Value = 5

Write-Host "done"
